# Auto-generated Excel COM-interop script applying the cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    # Force the cell to remain plain text (matches source data which stores
    # numeric-looking price strings like "42.424.70" as literal text, not numbers).
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "42.424.70"
Set-TextValue $ws.Range("E2") "  +0.50%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.285.27"
Set-TextValue $ws.Range("E3") "  -0.81%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.04%  "

# Row 5
Set-TextValue $ws.Range("D5") "322.85"
Set-TextValue $ws.Range("E5") "  +1.69%  "

# Row 6
Set-TextValue $ws.Range("D6") "102.40"
Set-TextValue $ws.Range("E6") "  -2.07%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.628"
Set-TextValue $ws.Range("E7") "  +0.26%  "

# Row 8
Set-TextValue $ws.Range("E8") "  +0.16%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.605"
Set-TextValue $ws.Range("E9") "  -0.48%  "

# Row 10
Set-TextValue $ws.Range("D10") "39.61"
Set-TextValue $ws.Range("E10") "  -0.32%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0904"
Set-TextValue $ws.Range("E11") "  -0.19%  "

# Row 12
Set-TextValue $ws.Range("D12") "8.31"
Set-TextValue $ws.Range("E12") "  -2.29%  "

# Row 13
Set-TextValue $ws.Range("E13") "  -0.47%  "

# Row 14
Set-TextValue $ws.Range("D14") "0.961"
Set-TextValue $ws.Range("E14") "  -1.17%  "

# Row 15
Set-TextValue $ws.Range("D15") "15.09"
Set-TextValue $ws.Range("E15") "  -2.03%  "

# Row 16
Set-TextValue $ws.Range("D16") "2.628.70"
Set-TextValue $ws.Range("E16") "  -0.93%  "

# Row 17
Set-TextValue $ws.Range("D17") "2.286.75"
Set-TextValue $ws.Range("E17") "  -1.01%  "

# Row 18
Set-TextValue $ws.Range("D18") "42.239.49"
Set-TextValue $ws.Range("E18") "  +0.38%  "

# Row 19
Set-TextValue $ws.Range("D19") "7.36"
Set-TextValue $ws.Range("E19") "  -4.86%  "

# Row 20
Set-TextValue $ws.Range("E20") "  -0.31%  "

# Row 21
Set-TextValue $ws.Range("D21") "12.73"
Set-TextValue $ws.Range("E21") "  +27.73%  "

# Row 22
Set-TextValue $ws.Range("B22") "PancakeSwap"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D22") "3.59"
Set-TextValue $ws.Range("E22") "  +0.28%  "

# Row 23
Set-TextValue $ws.Range("B23") "Litecoin"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D23") "72.93"
Set-TextValue $ws.Range("E23") "  -0.46%  "

# Row 24
Set-TextValue $ws.Range("D24") "267.49"
Set-TextValue $ws.Range("E24") "  -0.86%  "

# Row 25
Set-TextValue $ws.Range("E25") "  -3.48%  "

# Row 26
Set-TextValue $ws.Range("E26") "  -0.49%  "

# Row 27
Set-TextValue $ws.Range("D27") "10.84"
Set-TextValue $ws.Range("E27") "  -1.00%  "

# Row 28
Set-TextValue $ws.Range("D28") "2.32"
Set-TextValue $ws.Range("E28") "  +2.18%  "

# Row 29
Set-TextValue $ws.Range("D29") "22.45"
Set-TextValue $ws.Range("E29") "  -2.71%  "

# Row 30
Set-TextValue $ws.Range("D30") "38.04"
Set-TextValue $ws.Range("E30") "  +6.01%  "

# Row 31
Set-TextValue $ws.Range("D31") "164.34"
Set-TextValue $ws.Range("E31") "  -0.48%  "

# Row 32
Set-TextValue $ws.Range("D32") "6.04"
Set-TextValue $ws.Range("E32") "  +2.62%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.0873"
Set-TextValue $ws.Range("E33") "  -1.15%  "

# Row 34
Set-TextValue $ws.Range("E34") "  +1.17%  "

# Row 35
Set-TextValue $ws.Range("E35") "  -4.90%  "

# Row 36
Set-TextValue $ws.Range("D36") "2.54"
Set-TextValue $ws.Range("E36") "  -13.00%  "

# Row 37
Set-TextValue $ws.Range("E37") "  -1.30%  "

# Row 38
Set-TextValue $ws.Range("E38") "  +0.37%  "

# Row 39
Set-TextValue $ws.Range("D39") "3.68"
Set-TextValue $ws.Range("E39") "  +0.73%  "

# Row 40
Set-TextValue $ws.Range("D40") "2.75"
Set-TextValue $ws.Range("E40") "  -6.10%  "

# Row 41
Set-TextValue $ws.Range("E41") "  +2.41%  "

# Row 42
Set-TextValue $ws.Range("D42") "68.39"

# Row 43
Set-TextValue $ws.Range("B43") "Algorand"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D43") "0.225"
Set-TextValue $ws.Range("E43") "  -0.50%  "

# Row 44
Set-TextValue $ws.Range("B44") "FirstDigitalUSD"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D44") "1.00"
Set-TextValue $ws.Range("E44") "  -0.23%  "

# Row 45
Set-TextValue $ws.Range("D45") "90.14"
Set-TextValue $ws.Range("E45") "  -13.21%  "

# Row 46
Set-TextValue $ws.Range("D46") "12.16"
Set-TextValue $ws.Range("E46") "  +0.13%  "

# Row 47
Set-TextValue $ws.Range("D47") "113.38"
Set-TextValue $ws.Range("E47") "  -2.22%  "

# Row 48
Set-TextValue $ws.Range("D48") "80.42"
Set-TextValue $ws.Range("E48") "  +2.81%  "

# Row 49
Set-TextValue $ws.Range("E49") "  -1.48%  "

# Row 50
Set-TextValue $ws.Range("D50") "5.22"
Set-TextValue $ws.Range("E50") "  -2.04%  "

# Row 51
Set-TextValue $ws.Range("D51") "1.593.11"
Set-TextValue $ws.Range("E51") "  +2.42%  "
